$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the phone number used for test case 02 (row 2 / row 3 share it) ---
# Plain text, no leading-apostrophe needed in the final content.
$ws.Range("F2").Value = "'0823456789"
$ws.Range("F2").ClearFormats()

# --- Row 4 (test case 03) gets a *new* distinct text value that literally
#     starts with an apostrophe character. Doubling the leading apostrophe
#     forces text mode while keeping a single literal apostrophe in the text. ---
$ws.Range("F4").Value = "''0823456789"
$ws.Range("F4").ClearFormats()

# --- Append new row 7 (test case 06) ---
$ws.Range("A7").Value = "'06"
$ws.Range("A7").ClearFormats()

$ws.Range("B7").Value = "Ngoc Vu"

$ws.Range("C7").Value = "207 Giai Phong"

$ws.Range("D7").Value = "Ha Noi"

$ws.Range("E7").Value = "'100000"
$ws.Range("E7").ClearFormats()

$ws.Range("F7").Value = "'0823456789"
$ws.Range("F7").ClearFormats()

$ws.Range("G7").Value = "Hmm. We couldn’t find your device’s location. Please enter the address manually."

$ws.Range("H7").Value = "PASS"
